$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: remove the existing "_GoBack" bookmark that currently sits
# between "... Und damit " and " jede Familie ...", merging those two
# runs back into a single run (keeping the preceding single-letter "n"
# run intact/untouched). This must happen BEFORE we create the new
# "_GoBack" bookmark elsewhere, since Word bookmark names are unique —
# adding a second bookmark with the same name would just relocate this
# one instead of creating a new one.
# ---------------------------------------------------------------------

# Protect the seam before " die Besitzer der Katze ..." (right after the
# lone "n" run) with a temporary bookmark so the upcoming normalization
# edit cannot reach across it and swallow that run too.
$rngDie = $d.Content
$foundDie = $rngDie.Find.Execute(" die Besitzer der Katze")
if (-not $foundDie) {
    throw "Could not find ' die Besitzer der Katze' to protect"
}
$protectPoint = $d.Range($rngDie.Start, $rngDie.Start)
$d.Bookmarks.Add("TempProtect", $protectPoint)

# Drop the old "_GoBack" bookmark: this removes the barrier that was
# keeping " ... Und damit " and " jede Familie ..." apart.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Force Word's run-normalization by performing a genuine (round-tripped)
# text edit on one side of the old bookmark seam; this merges adjacent
# same-formatted runs up to the protecting bookmark, but does not cross
# it, so the "n" run is left alone.
$rngEdit = $d.Content
$foundEdit = $rngEdit.Find.Execute("Und damit")
if (-not $foundEdit) {
    throw "Could not find 'Und damit' to nudge for run normalization"
}
$rngEdit.Text = "Und damit!"

$rngRevert = $d.Content
$foundRevert = $rngRevert.Find.Execute("Und damit!")
if (-not $foundRevert) {
    throw "Could not find 'Und damit!' to revert"
}
$rngRevert.Text = "Und damit"

# Clean up the temporary protective bookmark.
$tempBm = $d.Bookmarks("TempProtect")
$tempBm.Delete()

# ---------------------------------------------------------------------
# Change 2: split the "... zu entwerfen. Also Programmiersprache ..."
# run into two runs around a new "_GoBack" bookmark, and fix the typo
# "Also" -> "Als" in the process.
# ---------------------------------------------------------------------

# Locate "Also" and insert the bookmark right after "Als" (i.e. before
# the trailing "o" of "Also"). Word.Bookmarks.Add on a collapsed range
# naturally splits the run it lands in without disturbing neighboring
# runs or their rsids.
$rngAlso = $d.Content
$foundAlso = $rngAlso.Find.Execute("Also")
if (-not $foundAlso) {
    throw "Could not find 'Also' to split/bookmark"
}
$splitPoint = $d.Range($rngAlso.Start + 3, $rngAlso.Start + 3)
$d.Bookmarks.Add("_GoBack", $splitPoint)

# Now remove the stray "o" that is left right after the bookmark
# (turning "Also" into "Als"), leaving the new run starting with
# " Programmiersprache ...".
$rngO = $d.Content
$foundO = $rngO.Find.Execute("o Programmiersprache")
if (-not $foundO) {
    throw "Could not find the leftover 'o' before 'Programmiersprache'"
}
$oChar = $d.Range($rngO.Start, $rngO.Start + 1)
$oChar.Delete()
